$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "CB00000"
$ws.Range("B6").Value = "test123"

$ws.Range("B6").Select()
